$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

# --- Update "Upcoming Meetings" (column D) dates for existing groups ---

# Row 4: Tucson Python Meetup
$ws.Range("D4").Value = "3/20;"

# Row 5: Google Developer Group (gains date-style number format)
$ws.Range("D5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = "3/8;"

# Row 6: Research Bazaar Arizona (already date-formatted)
$ws.Range("D6").Value = "3/8;"

# Row 8: Tucson WordPress Meetup
$ws.Range("D8").Value = "3/22;"

# Row 9: Tucson .NET User Group (already date-formatted)
$ws.Range("D9").Value = "3/15;"

# Row 11: Tucson Functional Programmers - becomes an actual date value
$ws.Range("D11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = "3/14/2018"

# Row 13: Tucson UX Rally
$ws.Range("D13").Value = "3/15;"

# --- Add new prospective groups ---

# Row 15: GDG Cloud Tucson
$ws.Range("A15").Value = "GDG Cloud Tucson"
$ws.Range("B15").Value = "https://www.meetup.com/GDG-Cloud-Tucson/"
$ws.Range("D15").NumberFormat = "d-mmm"
$ws.Range("D15").Value = "3/13;"
$ws.Range("E15").Value = 24
$ws.Range("G15").Value = "Andrew Slattery; "

# Row 16: R-Ladies Tucson
$ws.Range("A16").Value = "R-Ladies Tucson"
$ws.Range("B16").Value = "https://www.meetup.com/rladies-tucson-az/"
$ws.Range("C16").Value = "[Summarized from meetup desc] Promote gender diversity in the R statistical computing community. All events are intended for women but men are welcome to attend as a guest."
$ws.Range("D16").Value = "2/15;"
$ws.Range("E16").Value = 85
$ws.Range("G16").Value = "Adriana Picoral;"

# --- Update selection to match saved view state ---
$ws.Range("D9").Select() | Out-Null
